# Update to latest 4.0
# - Toggle guaranteed-dispatch flags (rows: hard coal, nuclear, solar thermal,
#   geothermal) on the DPbES sheet.
# - Clear the stray direct-formatting on About!A4 ("None").
# - Leave the workbook with the DPbES sheet active/selected (B10:AE10).

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsDPbES = $wb.Worksheets.Item("DPbES")

# hard coal (row 2): was excluded from guaranteed dispatch, now included
$wsDPbES.Range("B2:AE2").Value = 1

# nuclear (row 6): was included in guaranteed dispatch, now excluded
$wsDPbES.Range("B6:AE6").Value = 0

# solar thermal (row 10): was included in guaranteed dispatch, now excluded
$wsDPbES.Range("B10:AE10").Value = 0

# geothermal (row 12): was excluded from guaranteed dispatch, now included
$wsDPbES.Range("B12:AE12").Value = 1

# Clear the leftover direct formatting on the "None" label so it matches the
# other unstyled notes in the column.
$wsAbout.Range("A4").ClearFormats()

# DPbES is the sheet left on screen, scrolled/selected to the newly edited row.
$wsDPbES.Activate()
$wsDPbES.Range("B10:AE10").Select()
